# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2210", "_new" -> "_FV2304", then turn the header range into
# a proper Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldNames = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newNames = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

# Columns A-J (1-10) carry the "_old" suffix -> rename to "_FV2210"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $suffixless = $oldNames[$i].Substring(0, $oldNames[$i].Length - 4)
    $cell.Value = "$suffixless" + "_FV2210"
}

# Column K (11) is "diff" and stays as-is.

# Columns L-U (12-21) carry the "_new" suffix -> rename to "_FV2304"
for ($i = 0; $i -lt $newNames.Count; $i++) {
    $col = $i + 12
    $cell = $ws.Cells.Item(1, $col)
    $suffixless = $newNames[$i].Substring(0, $newNames[$i].Length - 4)
    $cell.Value = "$suffixless" + "_FV2304"
}

# Turn the used range into an Excel Table ("Table1") with headers.
$usedRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $usedRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (row 1) via a frozen pane.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
